# Add new columns I (I0) and J (IF) to the worksheet, matching the
# existing header row's style, and fill in the data values for rows 2-27.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for the two new columns.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Apply the same formatting as the existing header cells (bold, centered,
# bordered) by copying the format from H1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for rows 2-27.
$iValues = @(5,7,6,7,7,1,1,1,1,2,1,1,7,1,1,1,1,1,1,1,1,1,1,1,1,1)
$jValues = @(5,7,7,8,8,2,6,5,6,7,4,5,8,5,5,5,4,6,4,5,5,4,3,3,2,1)

for ($r = 0; $r -lt $iValues.Length; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$r]
    $ws.Cells.Item($row, 10).Value = $jValues[$r]
}
